# feat: cambio de contraseña por usuario desde navbar
# Updates the "Administrador" user's password (cell D3) from 1234 to tero2050,
# and leaves the sheet selection on D4 (as happens after editing D3 and
# pressing Enter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the Administrador row's Clave (password) value.
$ws.Range("D3").Value = "tero2050"

# Move / record the active selection to D4 (next cell down after the edit).
$ws.Range("D4").Select()
